# Applies the "Added gunicorn in requirements" data update to Sheet1:
#   - appends three new data rows (9, 10, 11) describing extra "Casting" phase
#     inspection records for the same B3F element as row 2.
#
# $wb / $ws refer to the already-open workbook / active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the exact date/time number format already used by column L (s="2" /
# numFmtId 165) so the new date cells share the same style index instead of
# creating a near-duplicate number format. The literal escaped format code
# (rather than the value read back from NumberFormat) is required so the
# engine recognizes it as the very same format already in the workbook.
$dateFormat = "yyyy\-mm\-dd\ hh:mm:ss"

function Set-TextCell($ws, $addr, [string]$text) {
    # Plain text values (never interpreted as numbers/dates by Excel)
    $ws.Range($addr).Value = $text
}

function Set-NumberCell($ws, $addr, $number) {
    $ws.Range($addr).Value = $number
}

function Set-DateCell($ws, $addr, $serial, $fmt) {
    $ws.Range($addr).NumberFormat = $fmt
    $ws.Range($addr).Value = $serial
}

function Set-QuotedTextCell($ws, $addr, [string]$text) {
    # Forces Excel to store a number-looking string ("0", "50", "24.0", ...)
    # as literal text (quote-prefixed) instead of silently converting it to
    # a numeric value.
    $ws.Range($addr).Value = "'" + $text
}

# ---------------------------------------------------------------------
# Row 9 - new "Casting: OK" reading for element 3d0f5ea4-...-b190b926
# ---------------------------------------------------------------------
Set-TextCell   $ws "A9" "3d0f5ea4-1394-46d0-b0b1-ba0ea9af8379"
Set-TextCell   $ws "B9" "Pilastro in calcestruzzo - Rettangolare"
Set-TextCell   $ws "C9" "Pilastro"
Set-TextCell   $ws "D9" "n\a"
Set-TextCell   $ws "E9" "via Merezzate, Milano>E10>P1"
Set-TextCell   $ws "F9" "C25/30"
Set-TextCell   $ws "G9" "Casting: OK"
Set-NumberCell $ws "H9" 0
Set-NumberCell $ws "I9" 0
Set-NumberCell $ws "J9" 1
Set-NumberCell $ws "K9" 1
Set-DateCell   $ws "L9" 43211.22545138889 $dateFormat
Set-TextCell   $ws "M9" "Appaltatore 1"
Set-NumberCell $ws "N9" 50
Set-NumberCell $ws "O9" 112
Set-TextCell   $ws "P9" "Bassa"
Set-TextCell   $ws "Q9" "Casting"
Set-NumberCell $ws "R9" 26
Set-NumberCell $ws "S9" 20
Set-NumberCell $ws "T9" 20
Set-TextCell   $ws "U9" "2018-09-10 05:00:13.436102"
Set-TextCell   $ws "V9" "40e526d7-263a-4f74-b935-1359b190b926"

# ---------------------------------------------------------------------
# Row 10 - new "Casting: Bad" reading for the same element
# ---------------------------------------------------------------------
Set-TextCell   $ws "A10" "3d0f5ea4-1394-46d0-b0b1-ba0ea9af8379"
Set-TextCell   $ws "B10" "Pilastro in calcestruzzo - Rettangolare"
Set-TextCell   $ws "C10" "Pilastro"
Set-TextCell   $ws "D10" "n\a"
Set-TextCell   $ws "E10" "via Merezzate, Milano>E10>P1"
Set-TextCell   $ws "F10" "C25/30"
Set-TextCell   $ws "G10" "Casting: Bad"
Set-NumberCell $ws "H10" 0
Set-NumberCell $ws "I10" 0
Set-NumberCell $ws "J10" 1
Set-NumberCell $ws "K10" 1
Set-DateCell   $ws "L10" 43211.22545138889 $dateFormat
Set-TextCell   $ws "M10" "Appaltatore 1"
Set-NumberCell $ws "N10" 50
Set-NumberCell $ws "O10" 112
Set-TextCell   $ws "P10" "Bassa"
Set-TextCell   $ws "Q10" "Casting"
Set-NumberCell $ws "R10" 24
Set-NumberCell $ws "S10" 26
Set-NumberCell $ws "T10" 26
Set-TextCell   $ws "U10" "2018-09-10 05:07:29.746704"
Set-TextCell   $ws "V10" "40e526d7-263a-4f74-b935-1359b190b926"

# ---------------------------------------------------------------------
# Row 11 - malformed/raw export row where the numeric columns were written
# out as text instead of numbers
# ---------------------------------------------------------------------
Set-TextCell        $ws "A11" "3d0f5ea4-1394-46d0-b0b1-ba0ea9af8379"
Set-TextCell        $ws "B11" "Pilastro in calcestruzzo - Rettangolare"
Set-TextCell        $ws "C11" "Pilastro"
Set-TextCell        $ws "D11" "n\a"
Set-TextCell        $ws "E11" "via Merezzate, Milano>E10>P1"
Set-TextCell        $ws "F11" "C25/30"
Set-TextCell        $ws "G11" "Casting: Bad"
Set-QuotedTextCell  $ws "H11" "0"
Set-QuotedTextCell  $ws "I11" "0"
Set-QuotedTextCell  $ws "J11" "1"
Set-QuotedTextCell  $ws "K11" "1"
Set-TextCell        $ws "L11" "2018-04-21 05:24:39"
Set-TextCell        $ws "M11" "Appaltatore 1"
Set-QuotedTextCell  $ws "N11" "50"
Set-QuotedTextCell  $ws "O11" "112"
Set-TextCell        $ws "P11" "Bassa"
Set-TextCell        $ws "Q11" "Casting"
Set-QuotedTextCell  $ws "R11" "24.0"
Set-QuotedTextCell  $ws "S11" "23.0"
Set-QuotedTextCell  $ws "T11" "25.0"
Set-TextCell        $ws "U11" "2018-09-10 05:33:02.786673"
Set-TextCell        $ws "V11" "Test"
